# Agile Tech.pptx - apply edits described by the commit:
#   "Fixed mention of encryption algorithm - We use SHA-256 encryption now"
#
# Changes:
#  1. Refresh the cached "datetimeFigureOut" date field text (14-10-03 -> 10/4/2014)
#     on the slide master and every slide layout's Date Placeholder.
#  2. Slide 3 ("Project Breakdown" / Features): merge the 3 runs of the
#     "Trade, add, remove, and increase shifts" bullet back into a single run.
#  3. Slide 5 ("Project Breakdown" / Technologies): replace "MD5" with "SHA-256"
#     in the "- MD5 Encryption for Passwords" bullet.
#  4. Slide 10 (Conclusion): replace "of" with "on" in the
#     "... while not cutting down of functionality." sentence.

$p = $ppt.ActivePresentation

function Get-DateShape($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $shp = $container.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            return $shp
        }
    }
    return $null
}

# Paragraphs() reports the trailing paragraph-mark (CR) in .Text/.Length
# for every paragraph except the very last one in a text body, so trim it
# off before doing exact string comparisons.
function Get-ParaText($para) {
    return $para.Text.TrimEnd([char]13)
}

# --- 1. Update the cached date field text everywhere it appears ---------
$master = $p.SlideMaster

$masterDateShape = Get-DateShape($master)
if ($masterDateShape -ne $null) {
    $masterDateShape.TextFrame.TextRange.Text = "10/4/2014"
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $dateShape = Get-DateShape($layout)
    if ($dateShape -ne $null) {
        $dateShape.TextFrame.TextRange.Text = "10/4/2014"
    }
}

# --- 2. Slide 3: merge "Trade, add, remove, and increase shifts" runs ---
$slide3 = $p.Slides.Item(3)
for ($j = 1; $j -le $slide3.Shapes.Count; $j++) {
    $shape = $slide3.Shapes.Item($j)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $paras = $tr.Paragraphs()
        for ($i = 1; $i -le $paras.Count; $i++) {
            $para = $tr.Paragraphs($i, 1)
            if ((Get-ParaText $para) -eq "Trade, add, remove, and increase shifts") {
                $full = $para.Characters(1, $para.Length)
                $full.Text = "Trade, add, remove, and increase shifts"
            }
        }
    }
}

# --- 3. Slide 5: "- MD5 Encryption for Passwords" -> SHA-256 ------------
$slide5 = $p.Slides.Item(5)
for ($j = 1; $j -le $slide5.Shapes.Count; $j++) {
    $shape = $slide5.Shapes.Item($j)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $paras = $tr.Paragraphs()
        for ($i = 1; $i -le $paras.Count; $i++) {
            $para = $tr.Paragraphs($i, 1)
            if ((Get-ParaText $para) -eq "- MD5 Encryption for Passwords") {
                $sub = $para.Characters(3, 4)
                $sub.Text = "SHA-256 "
            }
        }
    }
}

# --- 4. Slide 10: "... cutting down of functionality." -> "on" ----------
$slide10 = $p.Slides.Item(10)
for ($j = 1; $j -le $slide10.Shapes.Count; $j++) {
    $shape = $slide10.Shapes.Item($j)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $paras = $tr.Paragraphs()
        for ($i = 1; $i -le $paras.Count; $i++) {
            $para = $tr.Paragraphs($i, 1)
            if ((Get-ParaText $para) -eq "Our project aim is simplicity so that as little time is spent using the system, while not cutting down of functionality.") {
                $sub = $para.Characters(104, 3)
                $sub.Text = "on "
            }
        }
    }
}

Write-Output "edit complete"
